# Apply updated odds values to the "Jogos da Semana" worksheet
# per the source diff (2025-05-06 FlashScore export).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("I3").Value = 2.35
$ws.Range("L3").Value = 1.57
$ws.Range("M3").Value = 2.25
$ws.Range("N3").Value = 2.88
$ws.Range("O3").Value = 1.4
$ws.Range("P3").Value = 1.62
$ws.Range("Q3").Value = 2.2
$ws.Range("V3").Value = 15
$ws.Range("AF3").Value = 9.5

# Row 5
$ws.Range("G5").Value = 4.5
$ws.Range("H5").Value = 4.75
$ws.Range("I5").Value = 1.55
$ws.Range("J5").Value = 1.02
$ws.Range("K5").Value = 21
$ws.Range("R5").Value = 1.53
$ws.Range("S5").Value = 2.38
$ws.Range("T5").Value = 19
$ws.Range("AF5").Value = 10
$ws.Range("AH5").Value = 13

# Row 6
$ws.Range("N6").Value = 1.73
$ws.Range("O6").Value = 2.08

# Row 9
$ws.Range("J9").Value = 1.07
$ws.Range("K9").Value = 9
$ws.Range("N9").Value = 2.1
$ws.Range("O9").Value = 1.7

# Row 14
$ws.Range("G14").Value = 2.5
$ws.Range("I14").Value = 3
$ws.Range("J14").Value = 1.11
$ws.Range("K14").Value = 6.5
$ws.Range("T14").Value = 6.5
$ws.Range("U14").Value = 11
$ws.Range("AE14").Value = 7

# Row 16
$ws.Range("J16").Value = 1.07
$ws.Range("K16").Value = 9
$ws.Range("N16").Value = 2.15
$ws.Range("O16").Value = 1.67

# Row 18
$ws.Range("L18").Value = 1.4
$ws.Range("M18").Value = 2.75
$ws.Range("N18").Value = 2.25
$ws.Range("O18").Value = 1.62

# Row 22
$ws.Range("G22").Value = 2.38

# Row 23
$ws.Range("N23").Value = 2.35
$ws.Range("O23").Value = 1.57

# Row 30
$ws.Range("G30").Value = 2.36
$ws.Range("H30").Value = 3.35
$ws.Range("I30").Value = 2.65
$ws.Range("L30").Value = 1.26
$ws.Range("M30").Value = 3.5
$ws.Range("N30").Value = 1.81
$ws.Range("O30").Value = 1.89
$ws.Range("R30").Value = 1.67
$ws.Range("S30").Value = 2.08
$ws.Range("T30").Value = 7.2
$ws.Range("U30").Value = 9.800000000000001
$ws.Range("V30").Value = 7.6
$ws.Range("Y30").Value = 23
$ws.Range("Z30").Value = 9
$ws.Range("AA30").Value = 5
$ws.Range("AC30").Value = 50
$ws.Range("AE30").Value = 7.6
$ws.Range("AI30").Value = 18
$ws.Range("AJ30").Value = 25

# Row 31
$ws.Range("G31").Value = 4.1
$ws.Range("H31").Value = 3.6
$ws.Range("I31").Value = 1.74
$ws.Range("L31").Value = 1.22
$ws.Range("M31").Value = 3.8
$ws.Range("N31").Value = 1.71
$ws.Range("O31").Value = 2.02
$ws.Range("R31").Value = 1.68
$ws.Range("S31").Value = 2.06
$ws.Range("U31").Value = 20
$ws.Range("X31").Value = 30
$ws.Range("Y31").Value = 35
$ws.Range("Z31").Value = 10
$ws.Range("AA31").Value = 5.6
$ws.Range("AC31").Value = 50
$ws.Range("AE31").Value = 6.4
$ws.Range("AF31").Value = 7
$ws.Range("AH31").Value = 11
$ws.Range("AI31").Value = 11
$ws.Range("AJ31").Value = 20

# Row 33
$ws.Range("G33").Value = 2.87
$ws.Range("I33").Value = 2.15
$ws.Range("L33").Value = 1.21
$ws.Range("M33").Value = 3.55
$ws.Range("N33").Value = 1.65
$ws.Range("O33").Value = 2
$ws.Range("R33").Value = 1.57
$ws.Range("S33").Value = 2.12
$ws.Range("T33").Value = 11.25
$ws.Range("U33").Value = 16
$ws.Range("V33").Value = 10.5
$ws.Range("W33").Value = 35
$ws.Range("X33").Value = 22
$ws.Range("Y33").Value = 27
$ws.Range("Z33").Value = 13.5
$ws.Range("AA33").Value = 7.3
$ws.Range("AB33").Value = 13
$ws.Range("AC33").Value = 45
$ws.Range("AD33").Value = 300
$ws.Range("AE33").Value = 9.5
$ws.Range("AF33").Value = 11.5
$ws.Range("AG33").Value = 9
$ws.Range("AH33").Value = 21
$ws.Range("AI33").Value = 16
$ws.Range("AJ33").Value = 23

# Row 34
$ws.Range("G34").Value = 3
$ws.Range("H34").Value = 3.2
$ws.Range("I34").Value = 2.25
$ws.Range("L34").Value = 1.4
$ws.Range("M34").Value = 2.52
$ws.Range("N34").Value = 2.15
$ws.Range("O34").Value = 1.55
$ws.Range("P34").Value = 1.47
$ws.Range("Q34").Value = 2.35
$ws.Range("R34").Value = 1.93
$ws.Range("S34").Value = 1.7
$ws.Range("T34").Value = 7.7
$ws.Range("U34").Value = 14
$ws.Range("V34").Value = 11.5
$ws.Range("W34").Value = 37
$ws.Range("X34").Value = 30
$ws.Range("Y34").Value = 45
$ws.Range("Z34").Value = 7.7
$ws.Range("AA34").Value = 6.3
$ws.Range("AB34").Value = 17.5
$ws.Range("AC34").Value = 100
$ws.Range("AD34").Value = 101
$ws.Range("AE34").Value = 6.4
$ws.Range("AF34").Value = 9.75
$ws.Range("AG34").Value = 9.5
$ws.Range("AH34").Value = 21
$ws.Range("AI34").Value = 21
$ws.Range("AJ34").Value = 37

# Row 35
$ws.Range("G35").Value = 1.7
$ws.Range("H35").Value = 3.55
$ws.Range("I35").Value = 4.55
$ws.Range("L35").Value = 1.28
$ws.Range("M35").Value = 3.05
$ws.Range("N35").Value = 1.82
$ws.Range("O35").Value = 1.78
$ws.Range("P35").Value = 1.39
$ws.Range("Q35").Value = 2.55
$ws.Range("R35").Value = 1.78
$ws.Range("S35").Value = 1.83
$ws.Range("T35").Value = 6.8
$ws.Range("U35").Value = 8
$ws.Range("V35").Value = 8
$ws.Range("W35").Value = 13.5
$ws.Range("X35").Value = 13.5
$ws.Range("Y35").Value = 26
$ws.Range("Z35").Value = 10
$ws.Range("AA35").Value = 6.9
$ws.Range("AB35").Value = 15.5
$ws.Range("AC35").Value = 75
$ws.Range("AD35").Value = 600
$ws.Range("AE35").Value = 12.5
$ws.Range("AF35").Value = 26
$ws.Range("AG35").Value = 15
$ws.Range("AH35").Value = 80
$ws.Range("AI35").Value = 45
$ws.Range("AJ35").Value = 50

# Row 39
$ws.Range("J39").Value = 1.06
$ws.Range("K39").Value = 10
$ws.Range("N39").Value = 2.1
$ws.Range("O39").Value = 1.7

# Row 42
$ws.Range("H42").Value = 3.4
$ws.Range("N42").Value = 1.85
$ws.Range("O42").Value = 2
$ws.Range("P42").Value = 1.36
$ws.Range("Q42").Value = 3
$ws.Range("R42").Value = 1.67
$ws.Range("S42").Value = 2.1
$ws.Range("T42").Value = 9
$ws.Range("U42").Value = 12
$ws.Range("X42").Value = 17
$ws.Range("Z42").Value = 11
$ws.Range("AD42").Value = 151
$ws.Range("AE42").Value = 11
$ws.Range("AF42").Value = 17
$ws.Range("AJ42").Value = 29
